$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(7).Delete() | Out-Null
$ws.Range("B18").Select() | Out-Null
